$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TwoxTwoCET-Scalar")
$ws.Range("D2").Value = 0.84841617519703094
$ws.Range("E2").Value = 0.70716515735515251
$ws.Range("F2").Value = 0.38751484805497977
$ws.Range("D3").Value = 1.150632509520737
$ws.Range("E3").Value = 1.2893059702343586
$ws.Range("F3").Value = 1.597262943621075
$ws.Range("G3").Value = 1.9601317042077893
$ws.Range("E4").Value = 1.0068914996667164
$ws.Range("F4").Value = 0.99916851890038483
$ws.Range("G4").Value = 0.97073779227082335
$ws.Range("D5").Value = 1.0268107803031492
$ws.Range("E5").Value = 1.016185505434424
$ws.Range("F5").Value = 1.0862618717772727
$ws.Range("G5").Value = 1.1690468695684744
$ws.Range("D6").Value = 0.87672423219236773
$ws.Range("E6").Value = 0.88681325620491103
$ws.Range("F6").Value = 0.82408969313849223
$ws.Range("G6").Value = 0.76012909543532692
$ws.Range("D8").Value = 1.0351629708854644
$ws.Range("E8").Value = 1.0678356276899736
$ws.Range("F8").Value = 1.1033933594907175
$ws.Range("G8").Value = 1.1648853507249886
$ws.Range("D9").Value = 0.97200721521014122
$ws.Range("E9").Value = 0.94594737164343634
$ws.Range("F9").Value = 0.85817103861968957
$ws.Range("G9").Value = 0.7765902338166587
$ws.Range("D10").Value = 200
$ws.Range("E10").Value = 200
$ws.Range("F10").Value = 200
$ws.Range("G10").Value = 200
$ws.Range("D11").Value = 79.999999999999986
$ws.Range("E11").Value = 83.789979418307439
$ws.Range("D12").Value = 19.999999999999996
$ws.Range("E12").Value = 15.953297102016929
$ws.Range("F12").Value = 12.465814866890906
$ws.Range("E13").Value = 34.553784903611771
$ws.Range("F13").Value = 39.310669882974416
$ws.Range("G13").Value = 44.380086793200228
$ws.Range("E14").Value = 75.119576065481496
$ws.Range("F14").Value = 69.269028504233304
$ws.Range("D15").Value = 38.517354222141421
$ws.Range("E15").Value = 37.194397849630228
$ws.Range("F15").Value = 34.400605286275194
$ws.Range("D16").Value = 61.530004412497654
$ws.Range("E16").Value = 62.980517254317519
$ws.Range("F16").Value = 66.345864158482755
$ws.Range("D17").Value = 58.508040660383678
$ws.Range("E17").Value = 57.160534034089849
$ws.Range("F17").Value = 54.261106485862491
$ws.Range("G17").Value = 51.016980025031629
$ws.Range("D18").Value = 41.539717156383801
$ws.Range("E18").Value = 43.017230886987093
$ws.Range("F18").Value = 46.510809524574043
$ws.Range("G18").Value = 51.016980025031643
$ws.Range("C19").Value = 200
$ws.Range("D19").Value = 200
$ws.Range("E19").Value = 200
$ws.Range("F19").Value = 200
$ws.Range("G19").Value = 200
$ws.Range("C20").Value = 102.02649481767209
$ws.Range("D20").Value = 102.02649482407605
$ws.Range("E20").Value = 103.09328779209312
$ws.Range("F20").Value = 96.442586712999145
$ws.Range("G20").Value = 89.613092074379438
$ws.Range("C21").Value = 108.6294774886228
$ws.Range("D21").Value = 108.62947748112251
$ws.Range("E21").Value = 107.3936305887708
$ws.Range("F21").Value = 115.56763302716089
$ws.Range("G21").Value = 125.29200080619496
$ws.Range("C22").Value = 200.71701865447653
$ws.Range("D22").Value = 200
$ws.Range("E22").Value = 200
$ws.Range("F22").Value = 200
$ws.Range("G22").Value = 200
$ws.Range("C23").Value = 1.0035850932723827
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 1
$ws.Range("F23").Value = 1
$ws.Range("G23").Value = 1
$ws.Range("D24").Value = 1.0268107803031492
$ws.Range("E24").Value = 1.016185505434424
$ws.Range("F24").Value = 1.0862618717772727
$ws.Range("G24").Value = 1.1690468695684744
$ws.Range("D25").Value = 0.87672423219236773
$ws.Range("E25").Value = 0.88681325620491103
$ws.Range("F25").Value = 0.82408969313849223
$ws.Range("G25").Value = 0.76012909543532692
$ws.Range("D27").Value = 1.0351629708854644
$ws.Range("E27").Value = 1.0678356276899736
$ws.Range("F27").Value = 1.1033933594907175
$ws.Range("G27").Value = 1.1648853507249886
$ws.Range("D28").Value = 0.97200721521014122
$ws.Range("E28").Value = 0.94594737164343634
$ws.Range("F28").Value = 0.85817103861968957
$ws.Range("G28").Value = 0.7765902338166587
$ws.Range("D29").Value = 200
$ws.Range("E29").Value = 200
$ws.Range("F29").Value = 200
$ws.Range("G29").Value = 200
